$d = $word.ActiveDocument

# --- Change 1: merge the three runs of the "vehicle routing" sentence into one run ---
$sentence = "Developed and deployed software applications for intelligent vehicle routing optimization and supply chain network design. "
$r1 = $d.Content
$null = $r1.Find.Execute($sentence, $true, $false, $false, $false, $false, $true, 1, $false, $sentence, 2)

# --- Change 2: split " on premise or in the cloud..." -> replace "or" with "and" as its own run ---
$r2 = $d.Content
$null = $r2.Find.Execute("premise or in", $true, $true, $false, $false, $false, $true, 1, $false, "", 0)
$start = $r2.Start
$snippet = "premise or in"
$orIdx = $snippet.IndexOf("or")
$orStart = $start + $orIdx
$orEnd = $orStart + 2
$orRange = $d.Range($orStart, $orEnd)
$orRange.Text = "and"
# Force a run split by touching formatting on the freshly-inserted text
$andRange = $d.Range($orStart, $orStart + 3)
$andRange.Font.Color = 255
$andRange.Font.Color = 0

# --- Change 3: delete the trailing empty paragraph at the end of the document ---
$count = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($count)
$prev = $d.Paragraphs.Item($count - 1)
$delRange = $d.Range($prev.Range.End - 1, $last.Range.End)
$delRange.Delete()
